$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = 'Cluster name'
$ws.Cells.Item(1, 2).Value = 'Active cases'

# Data rows (cluster name, active cases)
$data = @(
    @('3035 Campbell Place Aged Care Glen Waverley', 11),
    @('3364 Assisi Centre Aged Care Rosanna', 21),
    @('3622 Olivet Care Aged Care Services Ringwood', 13),
    @('3633 Lifeview Emerald Glades Aged Care Emerald', 13),
    @('3961 Heritage Care Water Gardens Aged Care Facility Sydenham', 26),
    @('4167 Royal Freemasons Centennial Lodge Wantirna South', 21),
    @('AG Industries Pty Ltd Factory Thomastown', 17),
    @('Aintree Primary School Aintree', 18),
    @('Australian Meat Group Abattoir Dandenong South', 14),
    @('Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh', 31),
    @('Baden Powell College Tarneit', 15),
    @('Bandiana Primary School Bandiana', 10),
    @('Covenant College Bell Post Hill', 14),
    @('Gladstone Park Secondary College', 14),
    @('Hamlyn Banks Primary School Hamlyn Heights', 10),
    @('Hamlyn Views School Hamlyn Heights', 10),
    @('Hazelwood North Primary School Hazelwood North', 25),
    @('Hippity Hop Childcare and KindergartenPakenham', 10),
    @('Islamic College of Melbourne Tarneit Oct Nov', 16),
    @('M.C. Herd Corio', 10),
    @('Master Poultry Group West Footscray', 13),
    @('Morwell Park Primary School Morwell', 41),
    @('Nido Early School Woodend', 11),
    @('Northern Bay College Goldsworthy 9-12 Campus Corio', 18),
    @('Northern Bay College Wexford Campus Corio', 53),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 24),
    @('Oakleigh South Primary School Oakleigh South', 16),
    @('Our Lady''s Catholic Primary School Wangaratta', 12),
    @('Rutherglen Motor Inn and Walkabout Motel Rutherglen', 18),
    @('Sirius College Ibrahim Dellal Campus Sunshine', 12),
    @('Smartie Pants Early Learning and Development Diamond Creek', 19),
    @('St Brendans Primary School Shepparton', 10),
    @('St Georges Road Primary School Shepparton', 14),
    @('St Joseph''s School Quarry Hill', 29),
    @('St Josephs Catholic Primary School Warragul', 12),
    @('St Louis de Montfort''s School Aspendale', 13),
    @('St Paul''s Primary School Sunshine West', 12),
    @('St Vincents Hospital Emergency DepartmentMelbourne', 13),
    @('Stockdale Road Primary School Traralgon', 32),
    @('Story House Early Learning Epping October', 12),
    @('Sunbury Primary School Sunbury', 11),
    @('TUROSI PTY LTD Thomastown', 11),
    @('Templestowe Park Primary School Templestowe', 29),
    @('The Lake Primary School Cabarita', 18),
    @('The Royal Children''s Hospital MelbourneEmergency Department Parkville Tier 1A', 10),
    @('Warragul Regional College Warragul', 11),
    @('Werribee Mercy Hospital Emergency Department', 21),
    @('Western Health Sunshine Hospital Emergency Department St Albans', 10),
    @('Wodonga Primary School Wodonga', 21),
    @('Wodonga Senior Secondary College Wodonga', 14),
    @('Wodonga South Primary School Wodonga', 19),
    @('Woodend Primary School Woodend', 16),
    @('Wyndham Christian College Wyndham Vale', 14),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

Write-Host "Updated $($data.Count) data rows plus header"